# Edit: insert two new weekly price rows at the top of the "Zanahoria" price
# history table (rows 386-387), pushing all existing records down by two rows.
#
# Before: data rows occupy 2..414 (dimension A1:R414)
# After:  data rows occupy 2..416 (dimension A1:R416), with two brand-new rows
#         (Primera / Segunda for fecha 45106) inserted right before the old
#         row 386, and everything that used to be on rows 386-414 now living
#         on rows 388-416 (values unchanged, only their row number shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 386 and 387; this shifts old rows 386-414 down to
# 388-416 automatically (values travel with their rows), and Excel takes care
# of updating the sheet dimension for us.
$ws.Range("A386:A387").EntireRow.Insert()

# New row 386: Zanahoria, Primera, Vega Monumental Concepción, Bíobío
$ws.Range("A386").Value = 11
$ws.Range("B386").Value = "Vega Monumental Concepción"
$ws.Range("C386").Value = "Bíobío"
$ws.Range("D386").Value = 45106
$ws.Range("E386").Value = 8
$ws.Range("F386").Value = 100114013
$ws.Range("G386").Value = "Zanahoria"
$ws.Range("H386").Value = "Sin especificar"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 600
$ws.Range("K386").Value = 6000
$ws.Range("L386").Value = 6500
$ws.Range("M386").Value = 6250
$ws.Range("N386").Value = "$/saco 20 kilos"
$ws.Range("O386").Value = "Región de Ñuble"
$ws.Range("P386").Value = 312
$ws.Range("Q386").Value = 20
$ws.Range("R386").Value = "Hortaliza"

# New row 387: Zanahoria, Segunda, Vega Monumental Concepción, Bíobío
$ws.Range("A387").Value = 11
$ws.Range("B387").Value = "Vega Monumental Concepción"
$ws.Range("C387").Value = "Bíobío"
$ws.Range("D387").Value = 45106
$ws.Range("E387").Value = 8
$ws.Range("F387").Value = 100114013
$ws.Range("G387").Value = "Zanahoria"
$ws.Range("H387").Value = "Sin especificar"
$ws.Range("I387").Value = "Segunda"
$ws.Range("J387").Value = 300
$ws.Range("K387").Value = 5500
$ws.Range("L387").Value = 5500
$ws.Range("M387").Value = 5500
$ws.Range("N387").Value = "$/saco 20 kilos"
$ws.Range("O387").Value = "Región de Ñuble"
$ws.Range("P387").Value = 275
$ws.Range("Q387").Value = 20
$ws.Range("R387").Value = "Hortaliza"
